$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted at row 320, pushing every
# subsequent record (old rows 320..448) down by one row (new rows 321..449).
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new observation.
$ws.Cells.Item(320, 1).Value2  = 8
$ws.Cells.Item(320, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(320, 3).Value2  = "Coquimbo"
$ws.Cells.Item(320, 4).Value2  = 45009
$ws.Cells.Item(320, 5).Value2  = 4
$ws.Cells.Item(320, 6).Value2  = 100112003
$ws.Cells.Item(320, 7).Value2  = "Ajo"
$ws.Cells.Item(320, 8).Value2  = "Chino"
$ws.Cells.Item(320, 9).Value2  = "Primera"
$ws.Cells.Item(320, 10).Value2 = 440
$ws.Cells.Item(320, 11).Value2 = 17000
$ws.Cells.Item(320, 12).Value2 = 18000
$ws.Cells.Item(320, 13).Value2 = 17500
$ws.Cells.Item(320, 14).Value2 = '$/caja 10 kilos'
$ws.Cells.Item(320, 15).Value2 = "China"
$ws.Cells.Item(320, 16).Value2 = 1750
$ws.Cells.Item(320, 17).Value2 = 10
$ws.Cells.Item(320, 18).Value2 = "Hortaliza"
